# ---------------------------------------------------------------------------
# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" right after "总计" (before "2022-Q2").
# 2. Populate it with the fund-holding table for that quarter.
# 3. Update the "总计" (summary) sheet: push all existing quarter rows down
#    by one and insert the new 2022-Q3 totals at the top (row 2), appending
#    a brand-new row for what used to be the last row (2020-Q4).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$total = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Part A: add the new "2022-Q3" sheet right after "总计"
# ---------------------------------------------------------------------------
$new = $wb.Worksheets.Add($null, $total)
$new.Name = "2022-Q3"

# Re-fetch the (now shifted) "2022-Q2" sheet to use as a formatting donor.
$q2 = $wb.Worksheets.Item(3)

# Force columns B,D,E,F,G (and the header row) to text storage so numeric
# looking values ("31.61", "012210", ...) are kept as strings, not coerced
# to numbers (which would also drop leading zeros / trailing zeros).
# Column H (rank) is a genuine number and must stay General/numeric, so it
# is handled separately from this text block.
$textRange = $new.Range("B1:G12")
$textRange.NumberFormat = "@"
$headerTextCell = $new.Range("H1")
$headerTextCell.NumberFormat = "@"

# -- header row --------------------------------------------------------
$new.Cells.Item(1,2).Value = "基金代码"
$new.Cells.Item(1,3).Value = "基金名称"
$new.Cells.Item(1,4).Value = "基金规模"
$new.Cells.Item(1,5).Value = "股票总仓位"
$new.Cells.Item(1,6).Value = "仓位占比"
$new.Cells.Item(1,7).Value = "持有市值(亿元)"
$new.Cells.Item(1,8).Value = "仓位排名"

# -- data rows -----------------------------------------------------------
$new.Cells.Item(2,1).Value = 0
$new.Cells.Item(2,2).Value = "162006"
$new.Cells.Item(2,3).Value = "长城久富核心成长混合（LOF）A"
$new.Cells.Item(2,4).Value = "31.61"
$new.Cells.Item(2,5).Value = "71.49"
$new.Cells.Item(2,6).Value = "4.60"
$new.Cells.Item(2,7).Value = "1.4541"
$new.Cells.Item(2,8).Value = 4

$new.Cells.Item(3,1).Value = 1
$new.Cells.Item(3,2).Value = "398021"
$new.Cells.Item(3,3).Value = "中海能源策略混合"
$new.Cells.Item(3,4).Value = "21.73"
$new.Cells.Item(3,5).Value = "88.03"
$new.Cells.Item(3,6).Value = "3.70"
$new.Cells.Item(3,7).Value = "0.8040"
$new.Cells.Item(3,8).Value = 10

$new.Cells.Item(4,1).Value = 2
$new.Cells.Item(4,2).Value = "012210"
$new.Cells.Item(4,3).Value = "申万菱信智能汽车股票A"
$new.Cells.Item(4,4).Value = "2.98"
$new.Cells.Item(4,5).Value = "93.85"
$new.Cells.Item(4,6).Value = "5.22"
$new.Cells.Item(4,7).Value = "0.1556"
$new.Cells.Item(4,8).Value = 6

$new.Cells.Item(5,1).Value = 3
$new.Cells.Item(5,2).Value = "398061"
$new.Cells.Item(5,3).Value = "中海消费混合"
$new.Cells.Item(5,4).Value = "3.91"
$new.Cells.Item(5,5).Value = "85.30"
$new.Cells.Item(5,6).Value = "3.79"
$new.Cells.Item(5,7).Value = "0.1482"
$new.Cells.Item(5,8).Value = 9

$new.Cells.Item(6,1).Value = 4
$new.Cells.Item(6,2).Value = "014311"
$new.Cells.Item(6,3).Value = "大成优质精选混合A"
$new.Cells.Item(6,4).Value = "4.10"
$new.Cells.Item(6,5).Value = "47.17"
$new.Cells.Item(6,6).Value = "2.24"
$new.Cells.Item(6,7).Value = "0.0918"
$new.Cells.Item(6,8).Value = 9

$new.Cells.Item(7,1).Value = 5
$new.Cells.Item(7,2).Value = "015383"
$new.Cells.Item(7,3).Value = "长城久富核心成长混合（LOF）C"
$new.Cells.Item(7,4).Value = "1.82"
$new.Cells.Item(7,5).Value = "71.49"
$new.Cells.Item(7,6).Value = "4.60"
$new.Cells.Item(7,7).Value = "0.0837"
$new.Cells.Item(7,8).Value = 4

$new.Cells.Item(8,1).Value = 6
$new.Cells.Item(8,2).Value = "000976"
$new.Cells.Item(8,3).Value = "长城新兴产业灵活配置混合"
$new.Cells.Item(8,4).Value = "1.53"
$new.Cells.Item(8,5).Value = "71.73"
$new.Cells.Item(8,6).Value = "4.62"
$new.Cells.Item(8,7).Value = "0.0707"
$new.Cells.Item(8,8).Value = 4

$new.Cells.Item(9,1).Value = 7
$new.Cells.Item(9,2).Value = "012211"
$new.Cells.Item(9,3).Value = "申万菱信智能汽车股票C"
$new.Cells.Item(9,4).Value = "1.04"
$new.Cells.Item(9,5).Value = "93.85"
$new.Cells.Item(9,6).Value = "5.22"
$new.Cells.Item(9,7).Value = "0.0543"
$new.Cells.Item(9,8).Value = 6

$new.Cells.Item(10,1).Value = 8
$new.Cells.Item(10,2).Value = "000166"
$new.Cells.Item(10,3).Value = "中海信息产业精选混合"
$new.Cells.Item(10,4).Value = "0.77"
$new.Cells.Item(10,5).Value = "89.31"
$new.Cells.Item(10,6).Value = "5.11"
$new.Cells.Item(10,7).Value = "0.0393"
$new.Cells.Item(10,8).Value = 3

$new.Cells.Item(11,1).Value = 9
$new.Cells.Item(11,2).Value = "159804"
$new.Cells.Item(11,3).Value = "国寿安保国证创业板中盘精选88ETF"
$new.Cells.Item(11,4).Value = "1.10"
$new.Cells.Item(11,5).Value = "98.91"
$new.Cells.Item(11,6).Value = "1.94"
$new.Cells.Item(11,7).Value = "0.0213"
$new.Cells.Item(11,8).Value = 6

$new.Cells.Item(12,1).Value = 10
$new.Cells.Item(12,2).Value = "014312"
$new.Cells.Item(12,3).Value = "大成优质精选混合C"
$new.Cells.Item(12,4).Value = "0.07"
$new.Cells.Item(12,5).Value = "47.17"
$new.Cells.Item(12,6).Value = "2.24"
$new.Cells.Item(12,7).Value = "0.0016"
$new.Cells.Item(12,8).Value = 9

# Drop back to the default ("Normal") style everywhere in the data block so
# the only surviving style indices are the ones we apply below (matches the
# other quarter sheets, where only the header row + index column are styled).
$textRange.Style = "Normal"
$headerTextCell.Style = "Normal"

# -- formatting: bold/bordered header row + index column, copied from the
#    pre-existing "2022-Q2" sheet so the style id matches exactly. ---------
$srcHeader = $q2.Range("B1:H1")
$srcHeader.Copy()
$dstHeader = $new.Range("B1:H1")
$dstHeader.PasteSpecial(-4122)

$srcIndex = $q2.Range("A2:A5")
$srcIndex.Copy()
$dstIndex = $new.Range("A2:A12")
$dstIndex.PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Part B: update the "总计" (summary) sheet
# ---------------------------------------------------------------------------

# Create the new last row (row 9, was 2020-Q4's old row 8) by first cloning
# the formatting of the existing last index cell (A8) down into A9.
$srcLastA = $total.Range("A8")
$srcLastA.Copy()
$dstLastA = $total.Range("A9")
$dstLastA.PasteSpecial(-4122)
$total.Cells.Item(9,1).Value = 7

# Shift the quarter/count/value columns (B,C,D) down by one row, working
# from the bottom up so we never overwrite a value before reading it.
for ($r = 8; $r -ge 2; $r--) {
    $rNext = $r + 1
    $bVal = $total.Cells.Item($r, 2).Value()
    $cVal = $total.Cells.Item($r, 3).Value()
    $dVal = $total.Cells.Item($r, 4).Value()
    $total.Cells.Item($rNext, 2).Value = $bVal
    $total.Cells.Item($rNext, 3).Value = $cVal
    $total.Cells.Item($rNext, 4).Value = $dVal
}

# Write the brand-new 2022-Q3 totals into row 2 (column A index there is
# already 0 and needs no change).
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 11
$total.Cells.Item(2,4).Value = 2.92
